$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New master-data rows appended below the existing 21 data rows (rows 2-21),
# continuing the regcntr_id/usr_id/machine_id sequence.
$newRows = @(
    @(10002, 110021, 10021),
    @(10003, 110022, 10022),
    @(10004, 110023, 10023),
    @(10005, 110024, 10024),
    @(10006, 110025, 10025),
    @(10007, 110026, 10026),
    @(10008, 110027, 10027),
    @(10009, 110028, 10028),
    @(10010, 110029, 10029)
)

$r = 22
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = "eng"
    $ws.Cells.Item($r, 5).Value = $true
    $ws.Cells.Item($r, 6).Value = "superadmin"
    $ws.Cells.Item($r, 7).Value = "now()"
    $ws.Cells.Item($r, 8).Value = "now()"
    $r = $r + 1
}

# Match the saved selection state (active cell F14).
$ws.Range("F14").Select()

# Match the page setup (portrait orientation) recorded on the sheet.
$ws.PageSetup.Orientation = 1
